$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 574.0909
$ws.Range("I107").Value = 549.2857
$ws.Range("K107").Value = 549.2857
$ws.Range("M107").Value = 1370.7143
$ws.Range("H121").Value = 2220.2
$ws.Range("J121").Value = 2220.2
$ws.Range("L121").Value = 6660.599999999999
$ws.Range("N121").Value = -10154.6
$ws.Range("H137").Value = 4133
$ws.Range("I137").Value = 1393.8
$ws.Range("J137").Value = 100005
$ws.Range("K137").Value = 4181.4
$ws.Range("L137").Value = 300015
$ws.Range("M137").Value = -1631.4
$ws.Range("N137").Value = -305115
$ws.Range("H138").Value = 340221.8
$ws.Range("I138").Value = 5383.4116
$ws.Range("J138").Value = 439217.53
$ws.Range("K138").Value = 16150.2348
$ws.Range("L138").Value = 1317652.59
$ws.Range("M138").Value = -11010.2348
$ws.Range("N138").Value = -1327932.59
$ws.Range("H141").Value = 739
$ws.Range("I141").Value = 753.2857
$ws.Range("K141").Value = 2259.8571
$ws.Range("M141").Value = 2920.1429

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6051.431
$ws.Range("I32").Value = 5406.8726
$ws.Range("K32").Value = 5406.8726
$ws.Range("M32").Value = -5119.8726
$ws.Range("H61").Value = 3598.2954
$ws.Range("I61").Value = 1567.9706
$ws.Range("K61").Value = 1567.9706
$ws.Range("M61").Value = -1355.9706
$ws.Range("H74").Value = 244834.12
$ws.Range("I74").Value = 350418
$ws.Range("J74").Value = 3499.5715
$ws.Range("K74").Value = 350418
$ws.Range("L74").Value = 3499.5715
$ws.Range("M74").Value = -349544
$ws.Range("N74").Value = -5247.5715
$ws.Range("H77").Value = 244834.12
$ws.Range("I77").Value = 350418
$ws.Range("J77").Value = 3499.5715
$ws.Range("K77").Value = 1752090
$ws.Range("L77").Value = 17497.8575
$ws.Range("M77").Value = -1747722
$ws.Range("N77").Value = -26233.8575
$ws.Range("H110").Value = 2278.5186
$ws.Range("I110").Value = 1310.7059
$ws.Range("J110").Value = 3923.8
$ws.Range("K110").Value = 1310.7059
$ws.Range("L110").Value = 3923.8
$ws.Range("M110").Value = 734.2941000000001
$ws.Range("N110").Value = -8013.8
$ws.Range("H132").Value = 2480.3572
$ws.Range("I132").Value = 2055.6562
$ws.Range("J132").Value = 3839.4
$ws.Range("K132").Value = 6166.9686
$ws.Range("L132").Value = 11518.2
$ws.Range("M132").Value = -3636.9686
$ws.Range("N132").Value = -16578.2
$ws.Range("H135").Value = 119478
$ws.Range("J135").Value = 119478
$ws.Range("L135").Value = 119478
$ws.Range("N135").Value = -129618
$ws.Range("H136").Value = 3598.2954
$ws.Range("I136").Value = 1567.9706
$ws.Range("K136").Value = 4703.9118
$ws.Range("M136").Value = -2153.9118

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 322.92307
$ws.Range("I80").Value = 149
$ws.Range("J80").Value = 375.1
$ws.Range("K80").Value = 149
$ws.Range("L80").Value = 375.1
$ws.Range("M80").Value = 849
$ws.Range("N80").Value = -2371.1
$ws.Range("H83").Value = 322.92307
$ws.Range("I83").Value = 149
$ws.Range("J83").Value = 375.1
$ws.Range("K83").Value = 745
$ws.Range("L83").Value = 1875.5
$ws.Range("M83").Value = 4247
$ws.Range("N83").Value = -11859.5
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 73785.2
$ws.Range("J9").Value = 73785.2
$ws.Range("L9").Value = 73785.2
$ws.Range("N9").Value = -74121.2
$ws.Range("H134").Value = 2426.2415
$ws.Range("I134").Value = 2129.2693
$ws.Range("K134").Value = 6387.8079
$ws.Range("M134").Value = -3852.8079
$ws.Range("H140").Value = 85383.84
$ws.Range("J140").Value = 86749.164
$ws.Range("L140").Value = 86749.164
$ws.Range("N140").Value = -97109.164

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 23750
$ws.Range("J20").Value = 23750
$ws.Range("L20").Value = 23750
$ws.Range("N20").Value = -24240
$ws.Range("H122").Value = 3848521.8
$ws.Range("I122").Value = 3848521.8
$ws.Range("K122").Value = 11545565.4
$ws.Range("M122").Value = -11543115.4
$ws.Range("H140").Value = 73821.06
$ws.Range("J140").Value = 73821.06
$ws.Range("L140").Value = 73821.06
$ws.Range("N140").Value = -84181.06

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4033
$ws.Range("I68").Value = 5249.5
$ws.Range("J68").Value = 1600
$ws.Range("K68").Value = 5249.5
$ws.Range("L68").Value = 1600
$ws.Range("M68").Value = -4500.5
$ws.Range("N68").Value = -3098
$ws.Range("H71").Value = 4033
$ws.Range("I71").Value = 5249.5
$ws.Range("J71").Value = 1600
$ws.Range("K71").Value = 26247.5
$ws.Range("L71").Value = 8000
$ws.Range("M71").Value = -22503.5
$ws.Range("N71").Value = -15488
$ws.Range("H82").Value = 3142.9412
$ws.Range("I82").Value = 3370.5833
$ws.Range("J82").Value = 2596.6
$ws.Range("K82").Value = 3370.5833
$ws.Range("L82").Value = 2596.6
$ws.Range("M82").Value = -3009.5833
$ws.Range("N82").Value = -3318.6
$ws.Range("H85").Value = 3142.9412
$ws.Range("I85").Value = 3370.5833
$ws.Range("J85").Value = 2596.6
$ws.Range("K85").Value = 3370.5833
$ws.Range("L85").Value = 2596.6
$ws.Range("M85").Value = -2122.5833
$ws.Range("N85").Value = -5092.6
$ws.Range("H132").Value = 4467.3716
$ws.Range("I132").Value = 2818.4
$ws.Range("K132").Value = 8455.200000000001
$ws.Range("M132").Value = -5925.200000000001
$ws.Range("H136").Value = 6238.2
$ws.Range("I136").Value = 4626.0713
$ws.Range("J136").Value = 9999.833000000001
$ws.Range("K136").Value = 13878.2139
$ws.Range("L136").Value = 29999.499
$ws.Range("M136").Value = -11328.2139
$ws.Range("N136").Value = -35099.499

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 8336465
$ws.Range("I122").Value = 2940.1924
$ws.Range("J122").Value = 62504376
$ws.Range("K122").Value = 8820.5772
$ws.Range("L122").Value = 187513128
$ws.Range("M122").Value = -6370.5772
$ws.Range("N122").Value = -187518028
$ws.Range("H132").Value = 3795.8333
$ws.Range("J132").Value = 3489.8
$ws.Range("L132").Value = 10469.4
$ws.Range("N132").Value = -15529.4
